$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "W" in D1
$ws.Range("D1").Value = "W"

# Fill new column D with incrementing values 10..17
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 11
$ws.Range("D4").Value = 12
$ws.Range("D5").Value = 13
$ws.Range("D6").Value = 14
$ws.Range("D7").Value = 15
$ws.Range("D8").Value = 16
$ws.Range("D9").Value = 17

# Update existing value in B6
$ws.Range("B6").Value = 11.5

# Update selection to match target state
$ws.Range("B6").Select()
